# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text/code field ("004" -> "001"). Force text storage (avoid
# auto-conversion to the number 1) by temporarily applying a text number
# format, then clear the formatting again so no style index is left behind.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

# N2 is a text date/time stamp stored as a literal string.
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric financial figures for row 2.
$ws.Range("O2").Value = 7157222680.32
$ws.Range("P2").Value = 614690693.92
$ws.Range("Q2").Value = 510920179.46
$ws.Range("R2").Value = -27.0163305842
$ws.Range("S2").Value = 547173635.55
$ws.Range("T2").Value = -49.9040128234
$ws.Range("U2").Value = 39486096.72
$ws.Range("V2").Value = -93.2809131133
$ws.Range("W2").Value = 1344960189.92
$ws.Range("X2").Value = 570855444.1799999
$ws.Range("Y2").Value = -11.2053869966

# Z2 / AA2 were empty text cells and now hold numeric values.
$ws.Range("Z2").Value = 46263546.83
$ws.Range("AA2").Value = -21.3887687662

$ws.Range("AB2").Value = 5812262490.4
$ws.Range("AC2").Value = 1.6005581676
$ws.Range("AD2").Value = -3.4837567433
$ws.Range("AE2").Value = -20.6038125592
$ws.Range("AF2").Value = 301.9241023866
$ws.Range("AG2").Value = 18.7916493589
